# "Allow for multiple template documents for secondaries"
#
# The original "Offers" sheet has columns:
#   A=Quantity, B=Email, C=Address, D=PAN, E=Bank Account Number,
#   F=IFSC Code, G=Demat, H=City
#
# The edit inserts three new columns (First Name / Middle Name / Last Name)
# right after the Email column, relabels the first two header cells with
# "required field" style captions, and fills in sample first/last names for
# each of the four sample rows. Everything that used to live in columns
# C:H simply slides right to F:K.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: insert 3 blank columns before the old "Address" column
$ws.Range("C1:E1").EntireColumn.Insert()

# --- 2. Re-caption the columns that stayed in place
$ws.Range("A1").Value = "Offer Quantity *"
$ws.Range("B1").Value = "User (email) *"

# --- 3. Header captions for the three newly inserted columns
$ws.Range("C1").Value = "First Name *"
$ws.Range("D1").Value = "Middle Name"
$ws.Range("E1").Value = "Last Name *"

# The "Middle Name" / "Last Name *" headers pick up a distinct (no-fill)
# cell style, same as when a user explicitly clears the interior on a
# freshly typed header cell.
$ws.Range("D1:E1").Interior.ColorIndex = -4142

# Give the new columns the same display width as the neighbouring
# "Email" column (best effort - COM column widths round to pixels).
$ws.Range("C:E").ColumnWidth = $ws.Range("B1").ColumnWidth

# --- 4. Sample first / (blank middle) / last names for the four sample rows
$ws.Range("C2").Value = "Emp1"
$ws.Range("E2").Value = "John"

$ws.Range("C3").Value = "Emp2"
$ws.Range("E3").Value = "James"

$ws.Range("C4").Value = "Emp3"
$ws.Range("E4").Value = "Jim"

$ws.Range("C5").Value = "Emp4"
$ws.Range("E5").Value = "Jack"

# --- 5. Leave the selection where the author ended up (far right column)
$ws.Range("K1").Select()
